$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 data values (columns B:AH) to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 22.93
$ws.Range("C5").Value = 17.16
$ws.Range("D5").Value = 1.29
$ws.Range("E5").Value = 50.17
$ws.Range("F5").Value = 40.96
$ws.Range("G5").Value = 17.72
$ws.Range("H5").Value = 69.92
$ws.Range("I5").Value = 27.86
$ws.Range("J5").Value = 12.59
$ws.Range("K5").Value = 18.2
$ws.Range("L5").Value = 20.11
$ws.Range("M5").Value = 21.42
$ws.Range("N5").Value = 5.96
$ws.Range("O5").Value = 18.05
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 15.25
$ws.Range("R5").Value = 0.52
$ws.Range("S5").Value = 0.83
$ws.Range("T5").Value = 267.37
$ws.Range("U5").Value = 50.46
$ws.Range("V5").Value = 16.66
$ws.Range("W5").Value = 33.89
$ws.Range("X5").Value = 17.83
$ws.Range("Y5").Value = 2.35
$ws.Range("Z5").Value = 34.66
$ws.Range("AA5").Value = 14.71
$ws.Range("AB5").Value = 13.04
$ws.Range("AC5").Value = 15.31
$ws.Range("AD5").Value = 21.19
$ws.Range("AE5").Value = 0.48
$ws.Range("AF5").Value = 63.6
$ws.Range("AG5").Value = 9.35
$ws.Range("AH5").Value = 20.83

# Remove the last data row (row 6), reducing the used range to A1:AH5
$ws.Rows.Item(6).Delete()

# Narrow columns J (10) and V (22) from width 8 to width 7
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667
$ws.Columns.Item(22).ColumnWidth = 6.166666666666667
